$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.346.62"
$ws.Range("E2").Value = "'  +6.01%  "
$ws.Range("D3").Value = "'2.997.44"
$ws.Range("E3").Value = "'  +3.21%  "
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'580.42"
$ws.Range("E5").Value = "'  +2.38%  "
$ws.Range("D6").Value = "'162.69"
$ws.Range("E6").Value = "'  +12.45%  "
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("D8").Value = "'0.517"
$ws.Range("E8").Value = "'  +3.40%  "
$ws.Range("D9").Value = "'2.992.83"
$ws.Range("E9").Value = "'  +3.16%  "
$ws.Range("D10").Value = "'6.42"
$ws.Range("E10").Value = "'  -7.52%  "
$ws.Range("E11").Value = "'  +3.85%  "
$ws.Range("D12").Value = "'0.454"
$ws.Range("E12").Value = "'  +5.07%  "
$ws.Range("D13").Value = "'0.0000253"
$ws.Range("E13").Value = "'  +6.86%  "
$ws.Range("D14").Value = "'34.38"
$ws.Range("E14").Value = "'  +4.66%  "
$ws.Range("E15").Value = "'  -0.63%  "
$ws.Range("D16").Value = "'66.334.93"
$ws.Range("E16").Value = "'  +6.09%  "
$ws.Range("D17").Value = "'3.494.10"
$ws.Range("E17").Value = "'  +3.22%  "
$ws.Range("D18").Value = "'6.87"
$ws.Range("E18").Value = "'  +4.71%  "
$ws.Range("D19").Value = "'2.996.38"
$ws.Range("E19").Value = "'  +3.44%  "
$ws.Range("D20").Value = "'453.07"
$ws.Range("E20").Value = "'  +5.96%  "
$ws.Range("D21").Value = "'13.82"
$ws.Range("E21").Value = "'  +5.65%  "
$ws.Range("D22").Value = "'0.683"
$ws.Range("E22").Value = "'  +4.22%  "
$ws.Range("D23").Value = "'7.31"
$ws.Range("E23").Value = "'  +6.09%  "
$ws.Range("D24").Value = "'82.28"
$ws.Range("E24").Value = "'  +4.58%  "
$ws.Range("D25").Value = "'2.28"
$ws.Range("E25").Value = "'  +12.36%  "
$ws.Range("D26").Value = "'12.27"
$ws.Range("E26").Value = "'  +3.85%  "
$ws.Range("D27").Value = "'10.18"
$ws.Range("E27").Value = "'  +2.19%  "
$ws.Range("E28").Value = "'  -0.06%  "
$ws.Range("D29").Value = "'8.15"
$ws.Range("E29").Value = "'  +13.55%  "
$ws.Range("D30").Value = "'2.38"
$ws.Range("E30").Value = "'  +18.98%  "
$ws.Range("B31").Value = "'PEPE"
$ws.Range("C31").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "'0.0000104"
$ws.Range("E31").Value = "'  -4.69%  "
$ws.Range("B32").Value = "'PancakeSwap"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.62"
$ws.Range("E32").Value = "'  +5.59%  "
$ws.Range("D33").Value = "'27.13"
$ws.Range("E33").Value = "'  +5.00%  "
$ws.Range("E34").Value = "'  +4.60%  "
$ws.Range("E35").Value = "'  -0.07%  "
$ws.Range("D36").Value = "'0.988"
$ws.Range("E36").Value = "'  +4.37%  "
$ws.Range("D37").Value = "'5.78"
$ws.Range("E37").Value = "'  +7.28%  "
$ws.Range("D38").Value = "'2.08"
$ws.Range("E38").Value = "'  +8.70%  "
$ws.Range("D39").Value = "'49.51"
$ws.Range("E39").Value = "'  +1.76%  "
$ws.Range("D40").Value = "'2.93"
$ws.Range("E40").Value = "'  +0.74%  "
$ws.Range("D41").Value = "'0.304"
$ws.Range("E41").Value = "'  +13.20%  "
$ws.Range("B42").Value = "'Kaspa"
$ws.Range("C42").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.122"
$ws.Range("E42").Value = "'  +6.71%  "
$ws.Range("B43").Value = "'Arweave"
$ws.Range("C43").Value = "'https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").Value = "'43.88"
$ws.Range("E43").Value = "'  +6.09%  "
$ws.Range("D44").Value = "'8.42"
$ws.Range("E44").Value = "'  +4.22%  "
$ws.Range("D45").Value = "'399.93"
$ws.Range("E45").Value = "'  +10.56%  "
$ws.Range("D46").Value = "'0.0357"
$ws.Range("E46").Value = "'  +5.72%  "
$ws.Range("D47").Value = "'2.756.98"
$ws.Range("E47").Value = "'  +1.23%  "
$ws.Range("D48").Value = "'133.53"
$ws.Range("E48").Value = "'  -0.47%  "
$ws.Range("D50").Value = "'23.63"
$ws.Range("E50").Value = "'  +10.67%  "
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "'  +3.81%  "
